$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the old row 68 (shifts old rows 68:100 down to 71:103)
$ws.Rows("68:70").Insert()

# Populate the 3 newly inserted rows with this week's data (same market/region/category
# columns as the rest of the block; only the date + quality/volume/price columns change).
$newRows = @(
    @{ Row = 68; Fecha = 44609; Calidad = "Extra";   Volumen = 500; Min = 2800; Max = 2800; Prom = 2800; Unidad = "`$/unidad"; Origen = "Región de O'Higgins"; PrecioKg = 2800 },
    @{ Row = 69; Fecha = 44609; Calidad = "Primera"; Volumen = 500; Min = 2400; Max = 2400; Prom = 2400; Unidad = "`$/unidad"; Origen = "Región de O'Higgins"; PrecioKg = 2400 },
    @{ Row = 70; Fecha = 44609; Calidad = "Segunda"; Volumen = 500; Min = 2000; Max = 2000; Prom = 2000; Unidad = "`$/unidad"; Origen = "Región de O'Higgins"; PrecioKg = 2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = 100112028
    $ws.Cells.Item($row, 7).Value = "Sandia"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
